$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.222.53'
$ws.Range("E2").Value = '  +0.93%  '
$ws.Range("D3").Value = '1.689.02'
$ws.Range("E3").Value = '  +0.61%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.24'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.522'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.13'
$ws.Range("E8").Value = '  +13.54%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.262'
$ws.Range("E9").Value = '  +4.56%  '
$ws.Range("E10").Value = '  +1.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '1.927.48'
$ws.Range("E12").Value = '  +0.62%  '
$ws.Range("D13").Value = '1.701.72'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.559'
$ws.Range("E15").Value = '  +5.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.56'
$ws.Range("E16").Value = '  +2.81%  '
$ws.Range("D17").Value = '27.223.33'
$ws.Range("E17").Value = '  +0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.41'
$ws.Range("E18").Value = '  +0.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.12'
$ws.Range("E19").Value = '  -0.63%  '
$ws.Range("D20").Value = '0.0₃0745'
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("E21").Value = '  -0.03%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.57'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.65'
$ws.Range("E23").Value = '  +5.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.11'
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.67'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.32'
$ws.Range("E26").Value = '  +1.49%  '
$ws.Range("E27").Value = '  +2.53%  '
$ws.Range("E28").Value = '  +0.99%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0502'
$ws.Range("E30").Value = '  +0.65%  '
$ws.Range("E31").Value = '  +0.85%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.41'
$ws.Range("E32").Value = '  +2.40%  '
$ws.Range("D33").Value = '1.546.80'
$ws.Range("E33").Value = '  +4.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.24'
$ws.Range("E34").Value = '  +2.25%  '
$ws.Range("E36").Value = '  +4.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.605'
$ws.Range("E37").Value = '  +3.53%  '
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.06'
$ws.Range("E40").Value = '  +4.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.28'
$ws.Range("E41").Value = '  +3.00%  '
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").Value = '1.835.18'
$ws.Range("E45").Value = '  +0.94%  '
$ws.Range("E46").Value = '  +0.82%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '91.41'
$ws.Range("E47").Value = '  +0.86%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0110'
$ws.Range("E48").Value = '  +4.17%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.61'
$ws.Range("E49").Value = '  +5.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.27'
$ws.Range("E50").Value = '  +6.71%  '
$ws.Range("E51").Value = '  +1.82%  '
